# Insert a new "BKS" station row (alphabetically right after BDM, before BL67)
# into each of the three PSD sheets (HHE, HHN, HHZ). Inserting the row shifts
# every following station down by one row; the table's trailing duplicate
# "YBH" row falls off the bottom and is removed so the table keeps the same
# overall extent (rows 2-25).

$wb = $excel.ActiveWorkbook

$bksValues = @{
    "HHE" = @(-160, -131)
    "HHN" = @(-159, -130)
    "HHZ" = @(-178, -134)
}

foreach ($sheetName in @("HHE", "HHN", "HHZ")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Insert a blank row above the current row 3 (shifts BL67..YBH,YBH down by one)
    $ws.Rows.Item(3).Insert()

    # Populate the new row with the BKS station data
    $vals = $bksValues[$sheetName]
    $ws.Range("A3").Value = "BKS"
    $ws.Range("B3").Value = $vals[0]
    $ws.Range("C3").Value = $vals[1]

    # Remove the now-duplicated trailing YBH row (was row 25, shifted to row 26)
    $ws.Rows.Item(26).Delete()
}
